# Insert 3 new weekly price rows for Kiwi at Terminal Hortofrutícola Agro Chillán.
# The new rows are inserted right after the existing row 184, pushing all
# subsequent rows down by 3 (old row 185 -> new row 188, etc.) exactly like a
# native Excel "Insert Rows" operation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("185:187").Insert()

# Shared/fixed metadata for every row in this sheet (market + product identity).
$mercadoId   = 7
$mercado     = "Terminal Hortofrutícola Agro Chillán"
$region      = "Ñuble"
$codreg      = 16
$tipo        = "Fruta"
$productoId  = 100101
$producto    = "Berries"
$categoriaId = 100101007
$categoria   = "Kiwi"
$variedad    = "Hayward"
$unidad      = "`$/bandeja 18 kilos"
$kgUnidad    = 18

function Set-KiwiRow($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Origen, $PrecioKg) {
    $ws.Range("A$Row").Value = $mercadoId
    $ws.Range("B$Row").Value = $mercado
    $ws.Range("C$Row").Value = $region
    $ws.Range("D$Row").Value = $Fecha
    $ws.Range("E$Row").Value = $codreg
    $ws.Range("F$Row").Value = $tipo
    $ws.Range("G$Row").Value = $productoId
    $ws.Range("H$Row").Value = $producto
    $ws.Range("I$Row").Value = $categoriaId
    $ws.Range("J$Row").Value = $categoria
    $ws.Range("K$Row").Value = $variedad
    $ws.Range("L$Row").Value = $Calidad
    $ws.Range("M$Row").Value = $Volumen
    $ws.Range("N$Row").Value = $PrecioMin
    $ws.Range("O$Row").Value = $PrecioMax
    $ws.Range("P$Row").Value = $PrecioProm
    $ws.Range("Q$Row").Value = $unidad
    $ws.Range("R$Row").Value = $Origen
    $ws.Range("S$Row").Value = $PrecioKg
    $ws.Range("T$Row").Value = $kgUnidad
}

# New row 185 - Especial
Set-KiwiRow 185 45089 "Especial" 80 12000 12000 12000 "Región de O'Higgins" 667

# New row 186 - Primera
Set-KiwiRow 186 45089 "Primera" 80 10000 10000 10000 "Región de O'Higgins" 556

# New row 187 - Segunda
Set-KiwiRow 187 45089 "Segunda" 60 8000 8000 8000 "Región de O'Higgins" 444

Write-Host "Inserted 3 new rows (185-187); sheet now spans" $ws.UsedRange.Address()
